$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refactor: rename the header row (row 1) field/column names to match the
# updated ItemSO / WeaponData backing-field names (e.g. "id" -> "_id",
# "itemName" -> "_name", "icon" -> "_spritePath", "itemPrefab" -> "_prefabPath",
# "projectile" -> "_projectilePath", etc.)
$headers = @(
    "_id",
    "_name",
    "_description",
    "_spritePath",
    "_iconWidth",
    "_iconHeight",
    "_prefabPath",
    "_price",
    "_atk",
    "_atkRate",
    "_critRate",
    "_critDamage",
    "_range",
    "_lifeSteal",
    "_type",
    "_weaponTier",
    "_projectilePath",
    "_numberOfProjectile"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Move the visible selection/scroll position to reflect the new working area.
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("R1").Select()
